$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.465.96'
$ws.Cells.Item(2, 5).Value = '  -1.07%  '
$ws.Cells.Item(3, 4).Value = '2.523.65'
$ws.Cells.Item(3, 5).Value = '  -0.28%  '
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).Value = '''317.11'
$ws.Cells.Item(5, 5).Value = '  +3.77%  '
$ws.Cells.Item(6, 4).Value = '''94.47'
$ws.Cells.Item(6, 5).Value = '  -7.40%  '
$ws.Cells.Item(7, 5).Value = '  -0.79%  '
$ws.Cells.Item(9, 5).Value = '  -3.27%  '
$ws.Cells.Item(10, 4).Value = '''35.83'
$ws.Cells.Item(10, 5).Value = '  -4.87%  '
$ws.Cells.Item(11, 5).Value = '  -1.35%  '
$ws.Cells.Item(12, 4).Value = '''0.114'
$ws.Cells.Item(12, 5).Value = '  -0.16%  '
$ws.Cells.Item(13, 5).Value = '  -2.75%  '
$ws.Cells.Item(14, 4).Value = '2.911.29'
$ws.Cells.Item(14, 5).Value = '  -0.29%  '

# Rows 15-16: Chainlink / WrappedEther swap ranking positions with updated values
$ws.Cells.Item(15, 2).Value = 'WrappedEther'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(15, 4).Value = '2.531.74'
$ws.Cells.Item(15, 5).Value = '  -1.13%  '
$ws.Cells.Item(16, 2).Value = 'Chainlink'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(16, 4).Value = '''15.37'
$ws.Cells.Item(16, 5).Value = '  +1.08%  '

$ws.Cells.Item(17, 4).Value = '''0.845'
$ws.Cells.Item(17, 5).Value = '  -2.92%  '
$ws.Cells.Item(18, 4).Value = '42.515.88'
$ws.Cells.Item(18, 5).Value = '  -0.97%  '
$ws.Cells.Item(19, 4).Value = '''12.93'
$ws.Cells.Item(19, 5).Value = '  -2.03%  '
$ws.Cells.Item(21, 5).Value = '  -2.86%  '
$ws.Cells.Item(22, 4).Value = '''69.93'
$ws.Cells.Item(22, 5).Value = '  -2.33%  '
$ws.Cells.Item(23, 4).Value = '''250.38'
$ws.Cells.Item(23, 5).Value = '  -1.02%  '
$ws.Cells.Item(24, 4).Value = '''2.95'
$ws.Cells.Item(24, 5).Value = '  +0.70%  '
$ws.Cells.Item(25, 5).Value = '  -2.63%  '
$ws.Cells.Item(26, 4).Value = '''26.41'
$ws.Cells.Item(26, 5).Value = '  -2.84%  '
$ws.Cells.Item(27, 4).Value = '''0.997'
$ws.Cells.Item(27, 5).Value = '  -0.40%  '
$ws.Cells.Item(28, 4).Value = '''2.39'
$ws.Cells.Item(28, 5).Value = '  +2.83%  '
$ws.Cells.Item(29, 4).Value = '''10.15'
$ws.Cells.Item(29, 5).Value = '  -1.51%  '
$ws.Cells.Item(30, 5).Value = '  -0.39%  '
$ws.Cells.Item(31, 4).Value = '''5.98'
$ws.Cells.Item(31, 5).Value = '  -2.49%  '
$ws.Cells.Item(32, 4).Value = '''155.26'
$ws.Cells.Item(32, 5).Value = '  -1.38%  '
$ws.Cells.Item(33, 4).Value = '''19.16'
$ws.Cells.Item(33, 5).Value = '  +5.00%  '
$ws.Cells.Item(34, 5).Value = '  -0.66%  '
$ws.Cells.Item(35, 4).Value = '''3.26'
$ws.Cells.Item(35, 5).Value = '  -0.63%  '
$ws.Cells.Item(36, 5).Value = '  -1.71%  '
$ws.Cells.Item(37, 5).Value = '  -0.80%  '
$ws.Cells.Item(38, 5).Value = '  -4.45%  '
$ws.Cells.Item(39, 5).Value = '  -0.70%  '
$ws.Cells.Item(40, 4).Value = '''23.66'
$ws.Cells.Item(40, 5).Value = '  -1.11%  '
$ws.Cells.Item(41, 4).Value = '''2.33'
$ws.Cells.Item(41, 5).Value = '  +11.05%  '
$ws.Cells.Item(42, 5).Value = '  +0.28%  '
$ws.Cells.Item(43, 5).Value = '  -2.57%  '
$ws.Cells.Item(44, 5).Value = '  -1.63%  '
$ws.Cells.Item(45, 4).Value = '''3.24'
$ws.Cells.Item(45, 5).Value = '  -7.11%  '
$ws.Cells.Item(46, 4).Value = '2.015.50'
$ws.Cells.Item(46, 5).Value = '  -1.25%  '
$ws.Cells.Item(47, 4).Value = '''84.70'
$ws.Cells.Item(48, 4).Value = '''8.79'
$ws.Cells.Item(48, 5).Value = '  -1.66%  '
$ws.Cells.Item(49, 4).Value = '2.766.17'
$ws.Cells.Item(49, 5).Value = '  -0.43%  '
$ws.Cells.Item(50, 4).Value = '''74.06'
$ws.Cells.Item(50, 5).Value = '  +1.67%  '
$ws.Cells.Item(51, 4).Value = '''101.85'
$ws.Cells.Item(51, 5).Value = '  -1.08%  '
